$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "FUTBOL TALENTI"

$names = @(
    "Giovanni BOCCIA",
    "Michele CARUSO",
    "Mattia CASTELLUCCI",
    "Gabriele CAU",
    "Lorenzo CAU",
    "Gabriele CESANDRI",
    "Vincenzo CIRILLO",
    "Julian GENTILI",
    "Jacopo GENTILI",
    "Alessio IGLIOZZI",
    "Gianni LO PICCOLO",
    "Matteo  FERRAUTI",
    "Pietro MANCINI",
    "Flavio MOZZI",
    "Edoardo PALANCA",
    "Davide PENNACCHINI",
    "Daniele PICARIELLO",
    "Lorenzo RABBI",
    "Mathias NICOLETTI",
    "Davide RIFERZI"
)

$row = 22
foreach ($name in $names) {
    $ws.Cells.Item($row, 3).Value = $name
    $row = $row + 1
}

$ws.Range("C45").Value = "Edoardo BALESTRIERI, Lorenzo CAPANNOLO, Riccardo CASTELLANI, Enea COLANGELI, Gabriele GIANNOPOLO, Bruno MARRUCCI, Cristiano MIGENI, Lorenzo MIGENI, Tommaso BACIU, Valerio CAGNUCCI, Mirko BARONE, Eyad ELWANY"

$ws.Range("C52").Value = "Vincenzo GIANNOPOLI"
